# Insert a new row at row 57 (shifts existing rows 57..122 down to 58..123)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("57:57").Insert()

# Populate the newly inserted row 57 with the new data record.
# Columns A,B,C,E,F,G,H,I,N,Q,R copy the values that used to belong to
# (old) row 57, which is now row 58 after the insert; only D,J,K,L,M,O,P differ.
$ws.Range("A57").Value = 11
$ws.Range("B57").Value = "Vega Monumental Concepción"
$ws.Range("C57").Value = "Bíobío"
$ws.Range("D57").Value = 44930
$ws.Range("E57").Value = 8
$ws.Range("F57").Value = 100112024
$ws.Range("G57").Value = "Choclo"
$ws.Range("H57").Value = "Choclero"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 1800
$ws.Range("K57").Value = 250
$ws.Range("L57").Value = 300
$ws.Range("M57").Value = 272
$ws.Range("N57").Value = "`$/unidad"
$ws.Range("O57").Value = "Región Metropolitana"
$ws.Range("P57").Value = 272
$ws.Range("Q57").Value = 1
$ws.Range("R57").Value = "Hortaliza"
